$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data snapshot: rows 2-11 get re-shuffled with updated
# price/quality figures (row 8 unaffected).
$rows = @{
    2  = @{ D = 44489; L = "Primera"; M = 200; N = 24000; O = 25000; P = 24500; Q = "`$/caja 12 kilos";    S = 2042; T = 12 }
    3  = @{ D = 44167; L = "Segunda"; M = 200; N = 18000; O = 19000; P = 18500; Q = "`$/caja 13 kilos";    S = 1423; T = 13 }
    4  = @{ D = 44441; L = "Primera"; M = 100; N = 29000; O = 30000; P = 29500; Q = "`$/caja 12 kilos";    S = 2458; T = 12 }
    5  = @{ D = 44545; L = "Primera"; M = 200; N = 23000; O = 24000; P = 23500; Q = "`$/bandeja 12 kilos"; S = 1958; T = 12 }
    6  = @{ D = 44475; L = "Especial"; M = 200; N = 32000; O = 33000; P = 32500; Q = "`$/caja 12 kilos";   S = 2708; T = 12 }
    7  = @{ D = 44160; L = "Segunda"; M = 200; N = 19000; O = 20000; P = 19500; Q = "`$/caja 13 kilos";    S = 1500; T = 13 }
    9  = @{ D = 44468; L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/bandeja 10 kilos"; S = 2950; T = 10 }
    10 = @{ D = 44482; L = "Primera"; M = 160; N = 25000; O = 26000; P = 25500; Q = "`$/caja 12 kilos";    S = 2125; T = 12 }
    11 = @{ D = 44496; L = "Primera"; M = 200; N = 23000; O = 24000; P = 23500; Q = "`$/caja 12 kilos";    S = 1958; T = 12 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Range("D$r").Value = $data.D
    $ws.Range("L$r").Value = $data.L
    $ws.Range("M$r").Value = $data.M
    $ws.Range("N$r").Value = $data.N
    $ws.Range("O$r").Value = $data.O
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("S$r").Value = $data.S
    $ws.Range("T$r").Value = $data.T
}
